$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "order placement, fillo related changes" - refresh the Fillo test-data
# block (burhani/tno/vat/customer/company identifiers) used by the order
# placement scenario on this config sheet.
$ws.Range("B15").Value = "burhani1094"
$ws.Range("B16").Value = "tno400002317"
$ws.Range("B17").Value = "vat390002437"
$ws.Range("B20").Value = "customer3061"
$ws.Range("B21").Value = "CompanyName3029"

# Keep the selection where it was (B21) and nudge the view's scrolled
# column from A to B, matching the sheetView topLeftCell shift (A13 -> B13).
$ws.Range("B21").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 2
